# expansão das análises automáticas
# - E and F columns (particip / taxa_sucesso) were being stored as 0..1 fractions;
#   they now carry the already-multiplied-by-100 percentage number (display format
#   stays the same "0.00%" style already applied to those cells).
# - Three new computed columns are appended: apoio_medio, contribuicoes,
#   media_contribuicoes (L, M, N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rescale the "particip" (E) and "taxa_sucesso" (F) columns by 100.
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = 73.03370786516854
$ws.Range("F2").Value = 63.48717948717949

$ws.Range("E3").Value = 26.96629213483146
$ws.Range("F3").Value = 58.61111111111111

$ws.Range("E4").Value = 70.29972752043598
$ws.Range("F4").Value = 93.02325581395348

$ws.Range("E5").Value = 29.70027247956403
$ws.Range("F5").Value = 97.01834862385321

$ws.Range("E6").Value = 92.10526315789474
$ws.Range("F6").Value = 21.26984126984127

$ws.Range("E7").Value = 7.894736842105263
$ws.Range("F7").Value = 33.33333333333333

# ---------------------------------------------------------------------------
# 2) New header cells L1:N1 — copy K1's formatting (bold / centered / bordered
#    header style) then set the text.
# ---------------------------------------------------------------------------
$ws.Range("K1").Copy($ws.Range("L1"))
$ws.Range("K1").Copy($ws.Range("M1"))
$ws.Range("K1").Copy($ws.Range("N1"))

$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# ---------------------------------------------------------------------------
# 3) New data columns L, M, N for rows 2-7 (plain numbers, no special style —
#    matches the rest of the un-styled numeric columns C/D).
# ---------------------------------------------------------------------------
$ws.Range("L2").Value = 92.8307967260526
$ws.Range("M2").Value = 187904
$ws.Range("N2").Value = 303.5605815831987

$ws.Range("L3").Value = 87.50944228358395
$ws.Range("M3").Value = 75649
$ws.Range("N3").Value = 358.5260663507109

$ws.Range("L4").Value = 88.14455763556944
$ws.Range("M4").Value = 131371
$ws.Range("N4").Value = 136.8447916666667

$ws.Range("L5").Value = 93.84286760867113
$ws.Range("M5").Value = 72275
$ws.Range("N5").Value = 170.8628841607565

$ws.Range("L6").Value = 17.93451009587296
$ws.Range("M6").Value = 1926
$ws.Range("N6").Value = 14.37313432835821

$ws.Range("L7").Value = 30.65635216359388
$ws.Range("M7").Value = 282
$ws.Range("N7").Value = 15.66666666666667
